$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "12345678"
$ws.Range("C31").Style = "Normal"

$ws.Range("D31").Value = "asinha"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "19"
$ws.Range("E31").Style = "Normal"

$ws.Range("F31").Value = "QUALIDADE"
$ws.Range("G31").Value = "luis"
$ws.Range("H31").Value = "1º TURNO"
$ws.Range("I31").Value = "METRASCAN"
$ws.Range("J31").Value = "DISP SOLDA"
$ws.Range("K31").Value = "56ij67iumn"
$ws.Range("L31").Value = "C2025.0029"
